# Fixed update to excel issue
#
# 1. Rename the "Requested quantity" header on the existing sheets to more
#    descriptive names.
# 2. Add a new "PO Forecast" sheet (after "Monthly Trend") with a forecast
#    table (ds / PO_Forecast / yhat_lower / yhat_upper).

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- Rename headers on the existing sheets ---------------------------------
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet, positioned after "Monthly Trend" -----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "PO Forecast"

# Reuse the existing header style (bold, centered, bordered) and the
# existing date-serial number-format style from the Weekly Quantity sheet,
# so no new styles are created in the workbook.
$wsWeekly.Range("A1:B1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$ws.Range("A2:A19").PasteSpecial(-4122)

# --- Header row --------------------------------------------------------
$ws.Range("A1").Value = "ds"
$ws.Range("B1").Value = "PO_Forecast"
$ws.Range("C1").Value = "yhat_lower"
$ws.Range("D1").Value = "yhat_upper"

# --- Data rows -----------------------------------------------------------
$data = @(
    @(44934.99999999999, 105, 26.14963714840232, 176.6909544716457),
    @(44941.99999999999, 102, 23.15275049975291, 175.1296375519968),
    @(44962.99999999999, 94, 21.61384995795187, 165.196124962135),
    @(44969.99999999999, 91, 15.42788924926889, 164.4178572963136),
    @(44990.99999999999, 83, 6.57001962226484, 149.4670598457256),
    @(44997.99999999999, 80, 12.19043225580173, 152.0774562613223),
    @(45004.99999999999, 77, 3.109051249158437, 152.3001509243462),
    @(45011.99999999999, 74, 2.411920339565686, 141.9977227810361),
    @(45025.99999999999, 69, -2.710380067376696, 142.2136525871368),
    @(45032.99999999999, 66, -7.668767542907457, 142.6056969788045),
    @(45039.99999999999, 63, -5.219801807151994, 142.0427545330211),
    @(45046.99999999999, 61, -13.87951175949412, 137.3826405120001),
    @(45053.99999999999, 58, -15.94283844641383, 128.5815982752977),
    @(45060.99999999999, 55, -23.10265215820947, 129.3827597471412),
    @(45067.99999999999, 52, -22.37088343893523, 122.2635442215642),
    @(45074.99999999999, 50, -24.54276319635248, 128.0326083219512),
    @(45081.99999999999, 47, -33.50845657029253, 123.3047769805471),
    @(45088.99999999999, 44, -24.65405619488975, 114.2431912120348)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}
